$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the quote values in row 2
$ws.Range("F2").Value = 111
$ws.Range("G2").Value = 111.5

# Update the selected cell to H2
$ws.Range("H2").Select()
